$d = $word.ActiveDocument

# 1. Update the date in the first line (25.07.24 -> 24.07.24)
$d.Content.Find.Execute("25.07.24", $true, $false, $false, $false, $false,
                         $true, 1, $false, "24.07.24", 2)

# 2. Remove the Heading1 paragraph ("AI models collapse when trained on
#    recursively generated data") entirely.
$p = $d.Paragraphs.Item(2)
$p.Range.Delete()

# 3. Replace the paper title paragraph.
$d.Content.Find.Execute(
    "מאמר די חמוד שחוקר מה קורה שמאמנים מודלי AI על הדאטה הנוצר על ידי מודלי AI. בשתי מילים - לא הכל ורוד שם ויש כמה סיבות למה הדברים עלולים להשתבש:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The Empirical Impact of Neural Parameter Symmetries, or Lack Thereof", 2)

# 4. Replace the second paragraph's body text.
$d.Content.Find.Execute(
    "דאטה דריפט (איך זה בעברית?) קיצוני: אימון מודלים על דאטה שנוצרה על ידי מודלים אחרים גורם להתרחקות של התפלגות הדאטה הנוצר על ידי המודל החדש מהדאטה האמיתי (כלומר אגרגציה של מרחק בין ההתפלגויות שלהן)..",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "הסקירה היום תהיה קצרה וקלילה לעומת הסקירות האחרונות על מודלי דיפוזיה למיניהם. המאמר של היום חוקר סימטריות ברשתות נוירונים עמוקות. ניתן לראות די בקלות כי קיימות לא מעט פרמוטציות של המטרצות המשקלים בשכבות השונות של רשת שלמעשה לא משנות את המודל. כלומר אם תפעילו את המודל אחרי פרמוטציה על כל קלט תקבלו את אותה התוצאה כמו עם המודל המקורי.",
    2)

# 5. Replace the third paragraph's body text.
$d.Content.Find.Execute(
    "הבעיות מחמירות בזנבות התפלגות הדאטה (תחומים או שפות עם מעט דאטה למשל): ההידרדרות משפיעה בעיקר על זנבות התפלגות הדאטה, שם דאטה נדיר הופך להיות עוד פחות מיוצג",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "האם הסימטריות האלו מביאות לנו משהו טוב? בכלל לא בטול - לי זה נראה (למרות שאני לא מומחה גדול בתחום) כמו סוג של יתירות של יש במודלים שבלעדיה אולי ניתן היה להגיע למודלים קטנים יותר למשל. המאמר בוחן מה קורה במודל עם אנו מפרים את הסימטריה שיש במודל. אחת הדרכים להרוס את הסימטריה היא לקבע משקלות (לערכים אקראיים אך קבועים) במקומות שנבחרו באקראי במטריצות משקלים של הרשת. הדרך השניה היא להפעיל פונקציה אקטיבציה רק על המשקלים מסוימים.",
    2)

# 6. Replace the fourth paragraph's body text.
$d.Content.Find.Execute(
    "עוד יותר שגיאות: שגיאות בדאטה שנוצרו על ידי מודלים מצטברות לאורך דורות, מה שמוביל לירידה משמעותית בביצועים.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "המאמר חוקר איזה אפקטים מתרחשים אחרי שהורסים את הסימטריה במודל ומגלה כמה דברים די מעניינים….",
    2)

# 7. Remove the "קריסת השונות..." paragraph entirely. Find its index first
#    (it now sits right before the URL paragraph).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "קריסת השונות*") {
        $para.Range.Delete()
        break
    }
}

# 8. Replace the URL.
$d.Content.Find.Execute(
    "https://www.nature.com/articles/s41586-024-07566-y",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://arxiv.org/pdf/2405.20231",
    2)
